$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value into a cell while avoiding Excel auto-converting
# numeric-looking strings (e.g. "243.78") into actual numbers, and without
# leaving a residual number-format/quote-prefix style on the cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "42.400.32"
$ws.Range("E2").Value = "  -0.62%  "
Set-TextValue $ws.Range("D3") "2.238.67"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.19%  "
Set-TextValue $ws.Range("D5") "243.78"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("E6").Value = "  -0.41%  "
Set-TextValue $ws.Range("D7") "74.48"
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("E8").Value = "  +0.12%  "
Set-TextValue $ws.Range("D9") "0.613"
$ws.Range("E9").Value = "  -2.89%  "
Set-TextValue $ws.Range("D10") "42.94"
$ws.Range("E10").Value = "  -4.44%  "
Set-TextValue $ws.Range("D11") "0.0968"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("E12").Value = "  -4.65%  "
$ws.Range("E13").Value = "  +0.43%  "
Set-TextValue $ws.Range("D14") "2.574.47"
$ws.Range("E14").Value = "  -0.34%  "
Set-TextValue $ws.Range("D15") "14.39"
$ws.Range("E15").Value = "  -2.26%  "
Set-TextValue $ws.Range("D16") "0.844"
$ws.Range("E16").Value = "  -2.72%  "
Set-TextValue $ws.Range("D17") "2.292.33"
$ws.Range("E17").Value = "  +2.37%  "
Set-TextValue $ws.Range("D18") "42.251.55"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("E19").Value = "  +4.58%  "
Set-TextValue $ws.Range("D20") "6.25"
$ws.Range("E20").Value = "  +0.25%  "
Set-TextValue $ws.Range("D21") "73.24"
$ws.Range("E21").Value = "  +1.33%  "
Set-TextValue $ws.Range("D22") "11.15"
$ws.Range("E22").Value = "  -0.01%  "
Set-TextValue $ws.Range("D23") "231.76"
$ws.Range("E23").Value = "  -0.41%  "
Set-TextValue $ws.Range("D24") "2.10"
$ws.Range("E24").Value = "  -7.33%  "
$ws.Range("E25").Value = "  +0.14%  "
Set-TextValue $ws.Range("D26") "11.50"
$ws.Range("E26").Value = "  -3.49%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("E29").Value = "  -2.03%  "
Set-TextValue $ws.Range("D30") "167.40"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  -0.52%  "
Set-TextValue $ws.Range("D32") "5.77"
$ws.Range("E32").Value = "  +0.69%  "
Set-TextValue $ws.Range("D33") "0.0806"
$ws.Range("E33").Value = "  -2.54%  "
Set-TextValue $ws.Range("D34") "30.34"
$ws.Range("E34").Value = "  -7.27%  "
$ws.Range("E35").Value = "  -0.60%  "
Set-TextValue $ws.Range("D36") "0.109"
$ws.Range("E36").Value = "  -9.11%  "
Set-TextValue $ws.Range("D37") "4.39"
$ws.Range("E37").Value = "  -7.74%  "
$ws.Range("E38").Value = "  -4.22%  "
Set-TextValue $ws.Range("D39") "13.66"
$ws.Range("E39").Value = "  -5.04%  "
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("E41").Value = "  -1.65%  "
Set-TextValue $ws.Range("D42") "65.28"
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("E44").Value = "  -2.06%  "
Set-TextValue $ws.Range("D45") "105.30"
$ws.Range("E45").Value = "  -2.97%  "
Set-TextValue $ws.Range("D46") "0.101"
$ws.Range("E46").Value = "  -2.33%  "
Set-TextValue $ws.Range("D47") "2.38"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("E48").Value = "  -2.33%  "
Set-TextValue $ws.Range("D49") "1.17"
$ws.Range("E49").Value = "  -1.76%  "
Set-TextValue $ws.Range("D50") "2.68"
$ws.Range("E50").Value = "  -1.19%  "
Set-TextValue $ws.Range("D51") "2.447.67"
$ws.Range("E51").Value = "  -0.66%  "
